$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 279 (HIM-GENTLE BABY SOAP 75G)
$ws.Range("B279").Value = 48706
$ws.Range("E279").Value = 39.8
$ws.Range("F279").Value = -144
$ws.Range("G279").Value = -4795.2

# Row 280 (HIM-GENTLE BABY SOAP 75G)
$ws.Range("B280").Value = 64973
$ws.Range("E280").Value = 35.4
$ws.Range("F280").Value = 150
$ws.Range("G280").Value = 4995

# Row 313 (HUL-3Roses Dust [C] 500G Relaunch)
$ws.Range("B313").Value = 57854
$ws.Range("F313").Value = 2
$ws.Range("G313").Value = 611.6799999999999

# Row 314 (HUL-3Roses Dust [C] 500G Relaunch)
$ws.Range("B314").Value = 62997
$ws.Range("F314").Value = 72
$ws.Range("G314").Value = 22020.48

# Row 316 (HUL-Bru Inst Poly 50g)
$ws.Range("B316").Value = 63565
$ws.Range("E316").Value = 109.19
$ws.Range("F316").Value = 60
$ws.Range("G316").Value = 6162.6

# Row 317 (HUL-Bru Inst Poly 50g)
$ws.Range("B317").Value = 57077
$ws.Range("D317").Value = 93.08
$ws.Range("E317").Value = 111.2
$ws.Range("F317").Value = 1
$ws.Range("G317").Value = 93.08

# Row 318 (HUL-Bru Inst Poly 50g)
$ws.Range("B318").Value = 61610
$ws.Range("D318").Value = 102.71
$ws.Range("E318").Value = 122.71
$ws.Range("F318").Value = -58
$ws.Range("G318").Value = -5957.18

# Row 346 (HUL-Kissan nango jam 490g)
$ws.Range("B346").Value = 55373
$ws.Range("E346").Value = 163.62
$ws.Range("F346").Value = -94
$ws.Range("G346").Value = -13562.32

# Row 347 (HUL-Kissan nango jam 490g)
$ws.Range("B347").Value = 63520
$ws.Range("E347").Value = 153.4
$ws.Range("F347").Value = 97
$ws.Range("G347").Value = 13995.16

# Row 350 (HUL-Kissan Pineapple Jam 500G)
$ws.Range("B350").Value = 63571
$ws.Range("F350").Value = 29
$ws.Range("G350").Value = 4160.92

# Row 351 (HUL-Kissan Pineapple Jam 500G)
$ws.Range("B351").Value = 57802
$ws.Range("E351").Value = 162.71
$ws.Range("F351").Value = -79
$ws.Range("G351").Value = -11334.92

# Row 352 (HUL-Kissan Pineapple Jam 500G)
$ws.Range("B352").Value = 63531
$ws.Range("E352").Value = 152.53
$ws.Range("F352").Value = 80
$ws.Range("G352").Value = 11478.4

# Row 372 (HUL-Liril Soap 125 G)
$ws.Range("B372").Value = 57885
$ws.Range("E372").Value = 62.28
$ws.Range("F372").Value = 4
$ws.Range("G372").Value = 208.52

# Row 373 (HUL-Liril Soap 125 G)
$ws.Range("B373").Value = 63652
$ws.Range("E373").Value = 55.42
$ws.Range("F373").Value = 250
$ws.Range("G373").Value = 13032.5

# Row 379 (HUL-Lux Radiant Glow 3*150g)
$ws.Range("B379").Value = 63564
$ws.Range("E379").Value = 137.16
$ws.Range("F379").Value = 57
$ws.Range("G379").Value = 7353.57

# Row 380 (HUL-Lux Radiant Glow 3*150g)
$ws.Range("B380").Value = 61608
$ws.Range("E380").Value = 154.12
$ws.Range("F380").Value = -56
$ws.Range("G380").Value = -7224.56

# Row 382 (Hul-pears pure and gentle 3x125 gm)
$ws.Range("B382").Value = 63560
$ws.Range("E382").Value = 134.87
$ws.Range("F382").Value = 104
$ws.Range("G382").Value = 13193.44

# Row 383 (Hul-pears pure and gentle 3x125 gm)
$ws.Range("B383").Value = 60325
$ws.Range("E383").Value = 151.57
$ws.Range("F383").Value = -102
$ws.Range("G383").Value = -12939.72

# Row 389 (HUL-Rap Refresh Bolt 1Kg)
$ws.Range("B389").Value = 57817
$ws.Range("F389").Value = 3
$ws.Range("G389").Value = 239.43

# Row 390 (HUL-Rap Refresh Bolt 1Kg)
$ws.Range("B390").Value = 62865
$ws.Range("F390").Value = 151
$ws.Range("G390").Value = 12051.31

# Row 419 (HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp)
$ws.Range("B419").Value = 57856
$ws.Range("F419").Value = 2
$ws.Range("G419").Value = 342.66

# Row 420 (HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp)
$ws.Range("B420").Value = 63007
$ws.Range("F420").Value = 984
$ws.Range("G420").Value = 168588.72

# Row 421 (HUL-Surf Exl Mtc Liq Tl 1 Ltr Cp)
$ws.Range("B421").Value = 57857
$ws.Range("F421").Value = 3
$ws.Range("G421").Value = 453.51

# Row 422 (HUL-Surf Exl Mtc Liq Tl 1 Ltr Cp)
$ws.Range("B422").Value = 63008
$ws.Range("F422").Value = 504
$ws.Range("G422").Value = 76189.67999999999

# Row 431 (HUL-VIM BAR MULTIPACK FW 4X200G)
$ws.Range("B431").Value = 63102
$ws.Range("C431").Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Range("F431").Value = 36
$ws.Range("G431").Value = 2140.92

# Row 432 (HUL-Vim Bar Multipack Fw 4X200G)
$ws.Range("B432").Value = 53082
$ws.Range("C432").Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("F432").Value = 1
$ws.Range("G432").Value = 59.47

# Row 457 (JLM-MBD Shiny Toothbrush Safari)
$ws.Range("B457").Value = 31930
$ws.Range("E457").Value = 26.8
$ws.Range("F457").Value = -62
$ws.Range("G457").Value = -1390.04

# Row 458 (JLM-MBD Shiny Toothbrush Safari)
$ws.Range("B458").Value = 63681
$ws.Range("E458").Value = 23.84
$ws.Range("F458").Value = 65
$ws.Range("G458").Value = 1457.3

# Row 583 (CRE-Butter cremfills 100gm)
$ws.Range("B583").Value = 65066
$ws.Range("E583").Value = 13.61
$ws.Range("F583").Value = 313
$ws.Range("G583").Value = 4009.53

# Row 584 (CRE-Butter cremfills 100gm)
$ws.Range("B584").Value = 53263
$ws.Range("E584").Value = 15.29
$ws.Range("F584").Value = -309
$ws.Range("G584").Value = -3958.29

# Row 599 (CRE-Cremica Oatmeal Digestive 112.5 Gm)
$ws.Range("B599").Value = 64925
$ws.Range("E599").Value = 13.97
$ws.Range("F599").Value = 302
$ws.Range("G599").Value = 3971.3

# Row 600 (CRE-Cremica Oatmeal Digestive 112.5 Gm)
$ws.Range("B600").Value = 45709
$ws.Range("E600").Value = 15.69
$ws.Range("F600").Value = -300
$ws.Range("G600").Value = -3945

# Row 872 (Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm)
$ws.Range("B872").Value = 65079
$ws.Range("F872").Value = 21
$ws.Range("G872").Value = 858.27

# Row 873 (Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm)
$ws.Range("B873").Value = 65362
$ws.Range("F873").Value = 2
$ws.Range("G873").Value = 81.73999999999999
